$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 725. This pushes the existing row 725
# (and everything below it) down by two rows, creating two blank rows
# at 725 and 726 that we then populate with the new weekly data.
$ws.Rows("725:726").Insert()

# New row 725: "1a plateado" quality entry for the latest week.
$ws.Range("A725").Value = 4
$ws.Range("B725").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C725").Value = "Los Lagos"
$ws.Range("D725").Value = 45008
$ws.Range("E725").Value = 10
$ws.Range("F725").Value = "Fruta"
$ws.Range("G725").Value = 100102
$ws.Range("H725").Value = "Cítricos"
$ws.Range("I725").Value = 100102003
$ws.Range("J725").Value = "Limón"
$ws.Range("K725").Value = "Sin especificar"
$ws.Range("L725").Value = "1a plateado"
$ws.Range("M725").Value = 600
$ws.Range("N725").Value = 30000
$ws.Range("O725").Value = 31000
$ws.Range("P725").Value = 30500
$ws.Range("Q725").Value = "`$/malla 18 kilos"
$ws.Range("R725").Value = "Región de O'Higgins"
$ws.Range("S725").Value = 1694
$ws.Range("T725").Value = 18

# New row 726: "2a plateado" quality entry for the latest week.
$ws.Range("A726").Value = 4
$ws.Range("B726").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C726").Value = "Los Lagos"
$ws.Range("D726").Value = 45008
$ws.Range("E726").Value = 10
$ws.Range("F726").Value = "Fruta"
$ws.Range("G726").Value = 100102
$ws.Range("H726").Value = "Cítricos"
$ws.Range("I726").Value = 100102003
$ws.Range("J726").Value = "Limón"
$ws.Range("K726").Value = "Sin especificar"
$ws.Range("L726").Value = "2a plateado"
$ws.Range("M726").Value = 200
$ws.Range("N726").Value = 27000
$ws.Range("O726").Value = 27000
$ws.Range("P726").Value = 27000
$ws.Range("Q726").Value = "`$/malla 18 kilos"
$ws.Range("R726").Value = "Región de O'Higgins"
$ws.Range("S726").Value = 1500
$ws.Range("T726").Value = 18
